$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove "Felipe" row (row 5) entirely - no longer tracked.
$ws.Rows.Item(5).Delete()

# Remove the two oldest week columns (11_02_2024 in B, 18_02_2024 in C).
$ws.Range("B:C").EntireColumn.Delete()

# Add the new week column for 17_03_2024 with the days worked for that week.
$ws.Range("E1").Value = "17_03_2024"
$ws.Range("E2").Value = 5
$ws.Range("E3").Value = 5
$ws.Range("E4").Value = 3
$ws.Range("E5").Value = 3

$ws.Range("E5").Select()
